# feat: add search in browser
#
# Sheet1 gets two updates:
#   1. The "Data da publicacao" (publication date) cell L2, which was
#      previously blank, now holds 2023-01-01 (serial 44927). The cell
#      already carries the date-formatted style, so just set the value.
#   2. The view's scroll/selection state moves: the active cell/selection
#      becomes L3 (previously H5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Set the publication date value for the listing in row 2.
$ws.Range("L2").Value = "2023-01-01"

# 2) Move the selection to L3, which also updates the saved sheet view.
$ws.Range("L3").Select()
